$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "List der Charakterisierungsfaktoren" -> "Liste der Charakterisierungsfaktoren"
#    and move the "_GoBack" bookmark so it sits right after the newly
#    typed "e" (between "Liste" and " der Charakterisierungsfaktoren").
# -----------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("List der Charakterisierungsfaktoren")
if ($found) {
    $start = $r.Start

    # Type the missing "e" right after "List".
    $insPoint = $d.Range($start + 4, $start + 4)
    $insPoint.InsertAfter("e")

    # Nudge the freshly typed "e" so it keeps living in its own run
    # instead of being silently re-joined with the "List" run it
    # followed (round-tripping Bold leaves the visible formatting
    # untouched but keeps the run boundary).
    $eRange = $d.Range($start + 4, $start + 5)
    $eRange.Bold = 1
    $eRange.Bold = 0

    # Drop the (moved) "_GoBack" bookmark right after "Liste", before
    # the " der Charakterisierungsfaktoren" text.
    $bmRange = $d.Range($start + 5, $start + 5)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# -----------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that used to sit inside
#    "Schaltfl. Wirk. hinz." (a document can only have one).
# -----------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $goBack = $d.Bookmarks("_GoBack")
    $goBackStart = $goBack.Range.Start
    # Only the old location (inside "Wirk.") should be cleared here;
    # the one we just placed above is the new, correct one.
    if ($goBackStart -ne ($start + 5)) {
        $goBack.Delete()
    }
}
